# Auto-generated Excel COM-interop script
# Applies numeric corrections to the "Alpha_Profits" leve-profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 879
$ws.Range("I15").Value = 879
$ws.Range("K15").Value = 2637
$ws.Range("M15").Value = -2468
$ws.Range("H40").Value = 12733.333
$ws.Range("J40").Value = 5680
$ws.Range("L40").Value = 5680
$ws.Range("N40").Value = -6030
$ws.Range("H62").Value = 1851.25
$ws.Range("I62").Value = 1801.6666
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 1801.6666
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -1177.6666
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 1851.25
$ws.Range("I65").Value = 1801.6666
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 9008.333000000001
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -5888.333000000001
$ws.Range("N65").Value = -16240
$ws.Range("H98").Value = 3637.4666
$ws.Range("I98").Value = 4398.2
$ws.Range("J98").Value = 2116
$ws.Range("K98").Value = 4398.2
$ws.Range("L98").Value = 2116
$ws.Range("M98").Value = -2900.2
$ws.Range("N98").Value = -5112
$ws.Range("H105").Value = 79950
$ws.Range("J105").Value = 79950
$ws.Range("L105").Value = 79950
$ws.Range("N105").Value = -86938
$ws.Range("H106").Value = 3475
$ws.Range("I106").Value = 3475
$ws.Range("K106").Value = 3475
$ws.Range("M106").Value = -2844
$ws.Range("H122").Value = 3637.4666
$ws.Range("I122").Value = 4398.2
$ws.Range("J122").Value = 2116
$ws.Range("K122").Value = 13194.6
$ws.Range("L122").Value = 6348
$ws.Range("M122").Value = -10744.6
$ws.Range("N122").Value = -11248
$ws.Range("H132").Value = 1097.3695
$ws.Range("I132").Value = 1000.0238
$ws.Range("J132").Value = 2119.5
$ws.Range("K132").Value = 3000.0714
$ws.Range("L132").Value = 6358.5
$ws.Range("M132").Value = -470.0714000000003
$ws.Range("N132").Value = -11418.5
$ws.Range("H141").Value = 57465.11
$ws.Range("I141").Value = 60591.234
$ws.Range("K141").Value = 181773.702
$ws.Range("M141").Value = -176593.702

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5804.6816
$ws.Range("I32").Value = 5804.6816
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5804.6816
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -5517.6816
$ws.Range("N32").ClearContents()
$ws.Range("H74").Value = 1856.1786
$ws.Range("I74").Value = 2092.7334
$ws.Range("K74").Value = 2092.7334
$ws.Range("M74").Value = -1218.7334
$ws.Range("H77").Value = 1856.1786
$ws.Range("I77").Value = 2092.7334
$ws.Range("K77").Value = 10463.667
$ws.Range("M77").Value = -6095.667000000001
$ws.Range("H105").Value = 36739
$ws.Range("J105").Value = 36739
$ws.Range("L105").Value = 36739
$ws.Range("N105").Value = -43727
$ws.Range("H122").Value = 2149.5
$ws.Range("I122").Value = 1979.5
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 5938.5
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -3488.5
$ws.Range("N122").Value = -13898.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10091.444
$ws.Range("I86").Value = 2849.8
$ws.Range("J86").Value = 12876.692
$ws.Range("K86").Value = 2849.8
$ws.Range("L86").Value = 12876.692
$ws.Range("M86").Value = -1726.8
$ws.Range("N86").Value = -15122.692
$ws.Range("H89").Value = 10091.444
$ws.Range("I89").Value = 2849.8
$ws.Range("J89").Value = 12876.692
$ws.Range("K89").Value = 14249
$ws.Range("L89").Value = 64383.45999999999
$ws.Range("M89").Value = -8633
$ws.Range("N89").Value = -75615.45999999999
$ws.Range("H105").Value = 2386.6
$ws.Range("I105").Value = 2386.6
$ws.Range("K105").Value = 2386.6
$ws.Range("M105").Value = -639.5999999999999
$ws.Range("H134").Value = 2405.606
$ws.Range("I134").Value = 2603.6206
$ws.Range("K134").Value = 7810.861800000001
$ws.Range("M134").Value = -5275.861800000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1240.6364
$ws.Range("I99").Value = 1094.1111
$ws.Range("J99").Value = 1900
$ws.Range("K99").Value = 1094.1111
$ws.Range("L99").Value = 1900
$ws.Range("M99").Value = 403.8888999999999
$ws.Range("N99").Value = -4896
$ws.Range("H126").Value = 1240.6364
$ws.Range("I126").Value = 1094.1111
$ws.Range("J126").Value = 1900
$ws.Range("K126").Value = 3282.3333
$ws.Range("L126").Value = 5700
$ws.Range("M126").Value = -812.3333000000002
$ws.Range("N126").Value = -10640
$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040
$ws.Range("H134").Value = 2108.5278
$ws.Range("I134").Value = 2112.5938
$ws.Range("J134").Value = 2076
$ws.Range("K134").Value = 6337.7814
$ws.Range("L134").Value = 6228
$ws.Range("M134").Value = -3802.7814
$ws.Range("N134").Value = -11298

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I23").Value = 105
$ws.Range("J23").Value = 154.66667
$ws.Range("K23").Value = 315
$ws.Range("L23").Value = 464.00001
$ws.Range("M23").Value = -80
$ws.Range("N23").Value = -934.00001
$ws.Range("H122").Value = 2225.7778
$ws.Range("I122").Value = 1298.5
$ws.Range("K122").Value = 11686.5
$ws.Range("M122").Value = -9236.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2174.889
$ws.Range("I132").Value = 2047.8182
$ws.Range("J132").Value = 2734
$ws.Range("K132").Value = 6143.4546
$ws.Range("L132").Value = 8202
$ws.Range("M132").Value = -3613.4546
$ws.Range("N132").Value = -13262

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4156.9443
$ws.Range("I40").Value = 1868.5
$ws.Range("J40").Value = 7017.5
$ws.Range("K40").Value = 1868.5
$ws.Range("L40").Value = 7017.5
$ws.Range("M40").Value = -1732.5
$ws.Range("N40").Value = -7289.5
$ws.Range("H46").Value = 2739.303
$ws.Range("J46").Value = 3734.5454
$ws.Range("L46").Value = 3734.5454
$ws.Range("N46").Value = -4110.5454
$ws.Range("H68").Value = 3458.1177
$ws.Range("J68").Value = 2866
$ws.Range("L68").Value = 2866
$ws.Range("N68").Value = -4364
$ws.Range("H71").Value = 3458.1177
$ws.Range("J71").Value = 2866
$ws.Range("L71").Value = 14330
$ws.Range("N71").Value = -21818
$ws.Range("H132").Value = 8008.1
$ws.Range("I132").Value = 7016
$ws.Range("J132").Value = 9496.25
$ws.Range("K132").Value = 21048
$ws.Range("L132").Value = 28488.75
$ws.Range("M132").Value = -18518
$ws.Range("N132").Value = -33548.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4483.3335
$ws.Range("I62").Value = 4483.3335
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4483.3335
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3859.3335
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4483.3335
$ws.Range("I65").Value = 4483.3335
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 22416.6675
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -19296.6675
$ws.Range("N65").ClearContents()

